# Swap the data rows for "run 5 times each" (E=0.769, F=0.746) and
# "Stronger penalty + run 5 times each" (E=0.7, F=0.722) between row 7 and row 8.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current (pre-edit) contents of row 7 and row 8 for the columns
# that actually differ between the two rows (E, F, H - the G column is a
# formula that recalculates automatically from E/F).
$e7 = $ws.Range("E7").Value2
$f7 = $ws.Range("F7").Value2
$h7 = $ws.Range("H7").Value2

$e8 = $ws.Range("E8").Value2
$f8 = $ws.Range("F8").Value2
$h8 = $ws.Range("H8").Value2

# Write row 8's former values into row 7 ...
$ws.Range("E7").Value = $e8
$ws.Range("F7").Value = $f8
$ws.Range("H7").Value = $h8

# ... and row 7's former values into row 8.
$ws.Range("E8").Value = $e7
$ws.Range("F8").Value = $f7
$ws.Range("H8").Value = $h7

# Match the saved selection state (active cell moved from I7 to I9).
$ws.Range("I9").Select()

Write-Output "Swapped rows 7 and 8 (E, F, H columns) and updated selection to I9."
